$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column S = 19, Column T = 20 (1-based)
$colS = 19
$colT = 20

# Values to place in column T for rows 3..40 (as in the diff). $null means
# "copy format only, leave the cell empty" (matches the blank divider / group
# header rows that only carry a style in column S).
$values = @{
    3  = $null
    4  = 2021
    5  = $null
    6  = 1466
    7  = $null
    8  = 76
    9  = 15
    10 = 1
    11 = 188
    12 = 22
    13 = 15
    14 = "-"
    15 = "-"
    16 = 112
    17 = "-"
    18 = 6
    19 = "-"
    20 = 29
    21 = 1002
    22 = "-"
    23 = $null
    24 = 1029
    25 = $null
    26 = 51
    27 = 4
    28 = "-"
    29 = 127
    30 = 14
    31 = 12
    32 = "-"
    33 = "-"
    34 = 70
    35 = "-"
    36 = 3
    37 = "-"
    38 = 16
    39 = 732
    40 = "-"
}

# Rows where the new "total" figure needs a right-aligned variant of the
# existing wrap-text style (these introduce a brand-new cellXfs entry, same
# as the source workbook's diff).
$rightAlignRows = @(6, 24)

for ($r = 3; $r -le 40; $r++) {
    $src = $ws.Cells.Item($r, $colS)
    $dst = $ws.Cells.Item($r, $colT)

    # Copy the column S formatting into column T first so the new cell
    # inherits the same font/border/wrap style as its neighbour.
    $src.Copy()
    $dst.PasteSpecial(-4122)

    if ($rightAlignRows -contains $r) {
        $dst.HorizontalAlignment = -4152
    }

    $val = $values[$r]
    if ($null -ne $val) {
        $dst.Value = $val
    }
}

# Update the active selection to match the edited workbook (T3).
$ws.Range("T3").Select()

Write-Host "done"
